# Updated symbol list on Wed Dec 14 15:56:26 UTC 2022 with GitHub Actions
#
# This script applies the price/volume/coin updates described by the diff
# between the previous and new "cryptos.xlsx" snapshot. Numeric-looking
# values in column D are stored as text in the workbook (e.g. "272.09"),
# so we force the cell's number format to Text ("@") before assigning the
# value to avoid Excel silently re-interpreting the string as a float.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Column D price updates (rows 2-13) ---
Set-TextValue "D2"  "272.09"
Set-TextValue "D3"  "23.21"
Set-TextValue "D4"  "6.378"
Set-TextValue "D5"  "0.06288"
Set-TextValue "D6"  "3.651"
Set-TextValue "D7"  "6.739"
Set-TextValue "D8"  "1.387"
Set-TextValue "D9"  "0.8383"
Set-TextValue "D10" "0.1630"
Set-TextValue "D11" "0.08428"
Set-TextValue "D12" "0.03475"
Set-TextValue "D13" "0.03139"

# --- Rows 14 & 15 swap places: MCDex <-> BitMartToken ---
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09317"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.950"
$ws.Range("E15").Value = "14MCDexMCB"

# --- Remaining column D price updates ---
Set-TextValue "D16" "0.001721"
Set-TextValue "D17" "0.04835"
Set-TextValue "D18" "0.006269"
Set-TextValue "D22" "3.735"
Set-TextValue "D23" "2.329"
Set-TextValue "D24" "0.01392"
Set-TextValue "D27" "0.0003739"
Set-TextValue "D40" "0.04692"
Set-TextValue "D41" "0.006910"
Set-TextValue "D42" "0.1180"
Set-TextValue "D43" "0.003452"
Set-TextValue "D44" "0.01256"
Set-TextValue "D47" "0.7977"
Set-TextValue "D48" "0.09996"
Set-TextValue "D49" "0.00002098"
Set-TextValue "D50" "0.01239"

Write-Output "Applied cryptos.xlsx price/symbol updates"
